# Updates the Markov transition-matrix probabilities on Sheet1 to reflect
# the re-simulated game outcomes (more games simulated, optimization logic
# drafted). Each row of the matrix is a probability distribution over the
# 18 states (columns B:S) and still sums to 1 after the update.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    "B2"  = 0.2941176470588235
    "C2"  = 0.2941176470588235
    "J2"  = 0.08823529411764706
    "P2"  = 0.2647058823529412
    "S2"  = 0.05882352941176471

    "P4"  = 0.75
    "S4"  = 0.25

    "B6"  = 0.05555555555555555
    "D6"  = 0.05555555555555555
    "F6"  = 0.05555555555555555
    "O6"  = 0.05555555555555555
    "Q6"  = 0.05555555555555555
    "R6"  = 0.05555555555555555

    "B7"  = 0.1666666666666667
    "J7"  = 0.08333333333333333
    "O7"  = 0.08333333333333333
    "Q7"  = 0.25
    "R7"  = 0.08333333333333333

    "B8"  = 0.04761904761904762
    "F8"  = 0.1428571428571428
    "J8"  = 0.09523809523809523
    "O8"  = 0.09523809523809523
    "Q8"  = 0.09523809523809523
    "R8"  = 0.2380952380952381
    "S8"  = 0.2857142857142857

    "B10" = 0.1418439716312057
    "D10" = 0.02127659574468085
    "F10" = 0.06382978723404255
    "J10" = 0.1276595744680851
    "O10" = 0.02127659574468085
    "Q10" = 0.2340425531914894
    "R10" = 0.07801418439716312
    "S10" = 0.3120567375886525

    "G11" = 0.2592592592592592
    "J11" = 0.1111111111111111
    "K11" = 0.2962962962962963
    "L11" = 0.3333333333333333

    "G13" = 0.5
    "J13" = 0.5

    "H15" = 0.1578947368421053
    "J15" = 0.3157894736842105
    "M15" = 0.05263157894736842
    "O15" = 0.05263157894736842
    "S15" = 0.4210526315789473

    "F16" = 0.05
    "H16" = 0.1
    "I16" = 0.05
    "J16" = 0.65
    "K16" = 0.05
    "O16" = 0.05
    "S16" = 0.05

    "H17" = 0.05263157894736842
    "I17" = 0.02631578947368421
    "J17" = 0.631578947368421
    "K17" = 0.131578947368421
    "M17" = 0.02631578947368421
    "O17" = 0.05263157894736842
    "S17" = 0.07894736842105263

    "F18" = 0.05882352941176471
    "H18" = 0.2352941176470588
    "J18" = 0.5882352941176471
    "K18" = 0.05882352941176471
    "O18" = 0.05882352941176471

    "F19" = 0.01162790697674419
    "H19" = 0.1162790697674419
    "I19" = 0.01162790697674419
    "J19" = 0.5697674418604651
    "K19" = 0.1395348837209302
    "O19" = 0.05813953488372093
    "S19" = 0.09302325581395349
}

foreach ($cellRef in $updates.Keys) {
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
